# Apply the "added ifoCAST full series evaluation" update:
# each data row (quarters-ahead error series) shifts one column to the left
# (the oldest vintage/column is dropped) and a newly evaluated value is appended
# where a new vintage exists; rows that have no new vintage simply lose their
# last (now unused) cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowValues = @(-1.292459514723937, -1.353633693138139, 0.6878222440181159, 0.1052277069844749, 0.2009027477938327, -0.0595162892048901, 0.7759387069536011, 0.6758593470509796, 0.3220728409712834, 0.4702502314037287)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $rowValues[$i]
}

$rowValues = @(-3.501204164937531, -1.459748227781277, -2.042342764814918, -1.94666772400556, -2.207086761004283, -1.371631764845791, -1.471711124748413, -1.825497630828109, -1.677320240395664, -1.519874251391509)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(3, $i + 2).Value = $rowValues[$i]
}

$rowValues = @(1.980281758742053, 1.397687221708412, 1.49336226251777, 1.232943225519047, 2.068398221677539, 1.968318861774917, 1.614532355695221, 1.762709746127666, 1.920155735131821, 1.41949627860575)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(4, $i + 2).Value = $rowValues[$i]
}

$rowValues = @(1.458861400122613, 1.554536440931971, 1.294117403933248, 2.12957240009174, 2.029493040189118, 1.675706534109422, 1.823883924541867, 1.981329913546022, 1.480670457019951, 1.761028337555421)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $rowValues[$i]
}

$rowValues = @(-0.4869194962242832, -0.747338533223006, 0.0881164629354852, -0.01196289696713632, -0.3657494030468326, -0.2175720126143872, -0.06012602361023223, -0.560785480136303, -0.2804275996008339, -0.4524024574458455)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(6, $i + 2).Value = $rowValues[$i]
}

$rowValues = @(-0.164743996189365, 0.6707109999691262, 0.5706316400665047, 0.2168451339868084, 0.3650225244192538, 0.5224685134234088, 0.02180905689733798, 0.3021669374328071, 0.1301920795877955, 0.436054619334127)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(7, $i + 2).Value = $rowValues[$i]
}

$rowValues = @(0.5750359591597685, 0.4749565992571469, 0.1211700931774507, 0.269347483609896, 0.426793472614051, -0.07386598391201982, 0.2064918966234494, 0.0345170387784377, 0.3403795785247692, -0.271950374162517)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(8, $i + 2).Value = $rowValues[$i]
}

$rowValues = @(0.7353756362558697, 0.3815891301761735, 0.5297665206086188, 0.6872125096127738, 0.186553053086703, 0.4669109336221722, 0.2949360757771605, 0.600798615523492, -0.01153133716379418, 0.6768400480353174)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(9, $i + 2).Value = $rowValues[$i]
}

$rowValues = @(-0.4538658659823178, -0.3056884755498724, -0.1482424865457174, -0.6489019430717882, -0.3685440625363191, -0.5405189203813308, -0.2346563806349992, -0.8469863333222853, -0.1586149481231739, -0.4353061035472806)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(10, $i + 2).Value = $rowValues[$i]
}

$rowValues = @(-0.2056091156472509, -0.04816312664309591, -0.5488225831691667, -0.2684647026336975, -0.4404395604787092, -0.1345770207323777, -0.7469069734196638, -0.05853558822055238, -0.3352267436446591)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(11, $i + 2).Value = $rowValues[$i]
}
$ws.Cells.Item(11, 11).ClearContents()

$rowValues = @(0.3056233794366003, -0.1950360770894705, 0.08532180344599868, -0.08665305439901295, 0.2192094853473185, -0.3931204673399676, 0.2952509178591439, 0.01855976243503714)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(12, $i + 2).Value = $rowValues[$i]
}
$ws.Cells.Item(12, 10).ClearContents()

$rowValues = @(-0.3432134675219158, -0.06285558698644665, -0.2348304448314583, 0.0710320949148732, -0.541297857772413, 0.1470735274266985, -0.1296176279974082)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(13, $i + 2).Value = $rowValues[$i]
}
$ws.Cells.Item(13, 9).ClearContents()

$rowValues = @(-0.2203015759906016, -0.3922764338356133, -0.0864138940892818, -0.698743846776568, -0.01037246157745647, -0.2870636170015632)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(14, $i + 2).Value = $rowValues[$i]
}
$ws.Cells.Item(14, 8).ClearContents()

$rowValues = @(0.1083830226904575, 0.414245562436789, -0.1980843902504972, 0.4902869949486143, 0.2135958395245076)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(15, $i + 2).Value = $rowValues[$i]
}
$ws.Cells.Item(15, 7).ClearContents()

$rowValues = @(0.1338876819013198, -0.4784422707859664, 0.2099291144131452, -0.06676204101096155)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(16, $i + 2).Value = $rowValues[$i]
}
$ws.Cells.Item(16, 6).ClearContents()

$rowValues = @(-0.3064674129409547, 0.3819039722581568, 0.1052128168340501)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(17, $i + 2).Value = $rowValues[$i]
}
$ws.Cells.Item(17, 5).ClearContents()

$rowValues = @(0.07604143251182532, -0.2006497229122814)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(18, $i + 2).Value = $rowValues[$i]
}
$ws.Cells.Item(18, 4).ClearContents()

$rowValues = @(0.4116802297750048)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(19, $i + 2).Value = $rowValues[$i]
}
$ws.Cells.Item(19, 3).ClearContents()

$ws.Cells.Item(20, 2).ClearContents()
